# Refresh the cryptos list (prices / 1h volume %) per the GitHub Actions
# data-refresh commit. Column D (Price) holds values that look numeric
# (e.g. "65.677.14", "0.0780") but must stay plain text, exactly as they
# were authored as inline strings - so force text format before writing
# them, otherwise Excel's autodetect would coerce them into numbers/dates
# and mangle the formatting (dropped digits, scientific notation, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.677.14"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.447.15"
$ws.Range("E3").Value = "  -3.91%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.96"
$ws.Range("E5").Value = "  -1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.88"
$ws.Range("E6").Value = "  -7.74%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.447.81"
$ws.Range("E7").Value = "  -3.90%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +0.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.57"
$ws.Range("E10").Value = "  -4.52%  "
$ws.Range("E11").Value = "  -9.97%  "
$ws.Range("E12").Value = "  -8.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.027.96"
$ws.Range("E13").Value = "  -4.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000182"
$ws.Range("E14").Value = "  -11.28%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.52"
$ws.Range("E15").Value = "  -10.40%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.453.93"
$ws.Range("E16").Value = "  -3.46%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.566.92"
$ws.Range("E17").Value = "  -1.32%  "
$ws.Range("E18").Value = "  -2.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.83"
$ws.Range("E19").Value = "  -11.12%  "
$ws.Range("E20").Value = "  -8.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.77"
$ws.Range("E21").Value = "  -7.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "394.42"
$ws.Range("E22").Value = "  -6.90%  "
$ws.Range("E23").Value = "  -10.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.42"
$ws.Range("E24").Value = "  -6.15%  "
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.589.55"
$ws.Range("E26").Value = "  -3.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000106"
$ws.Range("E27").Value = "  -11.65%  "
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  -10.20%  "
$ws.Range("E30").Value = "  -9.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.20"
$ws.Range("E31").Value = "  -12.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.453.05"
$ws.Range("E32").Value = "  -3.67%  "
$ws.Range("E33").Value = "  -0.03%  "
$ws.Range("E34").Value = "  -6.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.90"
$ws.Range("E35").Value = "  -8.49%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "173.44"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.21"
$ws.Range("E37").Value = "  -14.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.91"
$ws.Range("E38").Value = "  -10.62%  "
$ws.Range("E39").Value = "  -8.63%  "
$ws.Range("E40").Value = "  -13.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0779"
$ws.Range("E41").Value = "  -8.72%  "
$ws.Range("E42").Value = "  -7.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.60"
$ws.Range("E43").Value = "  -5.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  +0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.44"
$ws.Range("E45").Value = "  -14.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.63"
$ws.Range("E46").Value = "  -12.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.02"
$ws.Range("E47").Value = "  -3.97%  "
$ws.Range("E48").Value = "  -2.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.55"
$ws.Range("E49").Value = "  -8.15%  "
$ws.Range("E50").Value = "  -16.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.207.16"
$ws.Range("E51").Value = "  -7.75%  "
